# Tighten spacing and tweak SFDC instance box
#
# This script:
#  1. Updates the "datetimeFigureOut" date placeholder text on the slide
#     master and all 11 slide layouts from "10/6/22" to "2022-10-06".
#  2. Repositions/resizes a handful of shapes on slide 1 (the "Salesforce
#     instance" box + its header icon, and the "Patient" icon/label/
#     connector) per the target diff.

# --- helpers ---------------------------------------------------------
# PowerPoint COM exposes shape geometry in points, but the underlying
# storage (and the diff we're matching) is in EMU (1 pt = 12700 EMU).
# The interop layer keeps Left/Top/Width/Height as 32-bit floats, so a
# naive emu/12700.0 conversion can truncate to one EMU below the target
# once round-tripped through float32. A tiny epsilon nudges it back onto
# the correct integral EMU value on save.
function EmuToPt($emu) {
    return ($emu / 12700.0) + 0.00002
}

$p = $ppt.ActivePresentation

# --- 1. Date placeholder fields ---------------------------------------
$oldDate = "10/6/22"
$newDate = "2022-10-06"

function Fix-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

Fix-DatePlaceholder $p.SlideMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Fix-DatePlaceholder $layouts.Item($j).Shapes
}

# --- 2. Shape geometry tweaks on slide 1 ------------------------------
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)

    switch ($shp.Id) {
        144 {
            # "Rectangle 143" - the "Salesforce instance" box: shrinks
            # from the left, right edge stays put, height/top unchanged.
            $shp.Left   = EmuToPt 2313187
            $shp.Top    = EmuToPt 386138
            $shp.Width  = EmuToPt 3189668
            $shp.Height = EmuToPt 5443164
        }
        148 {
            # "Graphic 147" - building icon in the box header, realigned
            # with the new left edge of the box.
            $shp.Left   = EmuToPt 2313187
            $shp.Top    = EmuToPt 390831
            $shp.Width  = EmuToPt 381000
            $shp.Height = EmuToPt 381000
        }
        208 {
            # "Graphic 22" - the Patient icon, moved to the right.
            $shp.Left   = EmuToPt 1539276
            $shp.Top    = EmuToPt 1326177
            $shp.Width  = EmuToPt 469900
            $shp.Height = EmuToPt 469900
        }
        209 {
            # "TextBox 39" - the "Patient" label, moved to the right.
            $shp.Left   = EmuToPt 1411345
            $shp.Top    = EmuToPt 1788610
            $shp.Width  = EmuToPt 724056
            $shp.Height = EmuToPt 261610
        }
        210 {
            # "Straight Arrow Connector 209" - connector from the
            # Patient icon to the Salesforce Health Cloud icon; resized
            # to the new shorter, perfectly horizontal span.
            $shp.Left   = EmuToPt 2009176
            $shp.Top    = EmuToPt 1561127
            $shp.Width  = EmuToPt 561850
            $shp.Height = EmuToPt 0
        }
    }
}
